$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$null = $lastPara.Range.InsertParagraphBefore()
$target = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ListadoCursosServlet</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Descripción: Nos devuelve el listado de cursos pertinentes dependiendo de si el usuario es Alumno o Profesor.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Parámetros entrada:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>request</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“profesor”: “profesor”}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Parámetros salida:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ok</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”: “ok”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>L</w:t></w:r><w:r><w:t>istCursos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>listCursos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>L</w:t></w:r><w:r><w:t>istGrupos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>listGrupos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”}, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>session</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “Token”: token, “role”: role}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ERROR</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”: “ERROR”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>errormsg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mensaje_error</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">”}, </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>session</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “Token”: token, “role”: role}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ListadoAlumnosServlet</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Descripción: Nos devuelve el listado de alumnos asignados a un curso y grupo especifico, uso disponible solo para el Profesor.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Parámetros entrada:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>request</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:r><w:t>G</w:t></w:r><w:r><w:t>rupo”: “grupo”}, {“</w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t>urso”: “curso”}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Parámetros salida:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ok</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”: “ok”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ListAlumnos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>list</w:t></w:r><w:r><w:t>alumnos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”}, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>session</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “Token”: token, “role”: role}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:lastRenderedPageBreak/><w:t>ERROR</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”: “ERROR”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>errormsg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mensaje_error</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">”}, </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>session</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “Token”: token, “role”: role}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>InformaciónAlumnoServlet</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Descripción: Nos muestra la información detallada de un alumno</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Parámetros entrada: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>request</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>idAlumno</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>emailAlumno</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”}, {“Curso”: “curso”</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>}</w:t></w:r><w:r><w:t>,{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>“Grupo”: “grupo”}</w:t></w:r><w:r><w:t>}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Parámetros salida:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ok</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”: “ok”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Informacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>informacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”}, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>session</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “Token”: token, “role”: role}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ERROR</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”: “ERROR”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>errormsg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mensaje_error</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">”}, </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>session</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, “Token”: token, “role”: role}}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$null = $target.InsertXML($xml)
